# Update testFile to include Time in DateTime (FsSpreadsheet test workbook).
#
# Changes applied to the "WithTable" sheet (the active sheet):
#   - C3: 2023-10-15 (date only)            -> 2023-10-15 18:00 (date + time)
#   - C4: 2023-10-16 (date only)            -> 2023-10-16 20:00 (date + time)
#   - Both cells get a custom "d/m/yy h:mm;@" number format (date + time).
#   - Column C is widened to fit the new, longer date/time display.
#   - The active cell / selection moves to E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WithTable")

# Give the DateTime values an explicit time-of-day component.
$ws.Range("C3").Value2 = 45214.75
$ws.Range("C4").Value2 = 45215.833333333336

# Apply a custom date+time number format to the updated cells so the time
# portion is visible (this allocates a new cellXf/numFmt pair).
$ws.Range("C3:C4").NumberFormat = "d/m/yy\ h:mm;@"

# Widen column C to accommodate the longer date/time text.
$ws.Columns("C:C").ColumnWidth = 25

# Move the selection/active cell.
$ws.Range("E11").Select()
